$wb = $excel.ActiveWorkbook

# Insert the new "Periods" sheet right after "Stages", before "Programs"
$stagesSheet = $wb.Worksheets.Item("Stages")
$periods = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $stagesSheet)
$periods.Name = "Periods"

# Header row
$periods.Range("A1").Value = "initialDate"
$periods.Range("B1").Value = "name"
$periods.Range("C1").Value = "programID"

# Data rows - values entered with a leading apostrophe so Excel stores them
# as text (quote-prefixed) even though they look like dates/numbers.
$periods.Range("A2").Value = "'2014-08-31"
# Apply the short-date display format to column A (the data still reads back
# as text because of the quote-prefix above). Format A2 first, then use the
# format painter (Copy + PasteSpecial formats) so A3:A4 share the exact same
# style record instead of Excel minting a new one per cell.
$periods.Range("A2").NumberFormat = "mm-dd-yy"
$periods.Range("B2").Value = "period1program1"
$periods.Range("C2").Value = "'31"

$periods.Range("A3").Value = "'2014-08-29"
$periods.Range("B3").Value = "period2program1"
$periods.Range("C3").Value = "'32"

$periods.Range("A4").Value = "'2014-08-28"
$periods.Range("B4").Value = "period3program1"
$periods.Range("C4").Value = "'33"

$periods.Range("A2").Copy() | Out-Null
$periods.Range("A3:A4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Selections / active cells
$periods.Range("A5").Select()

$programs = $wb.Worksheets.Item("Programs")
$programs.Range("F5").Select()

$stages = $wb.Worksheets.Item("Stages")
$stages.Range("B1").Select()

$periods.Activate()
